# Update metric data: Tue Apr 29 09:10:15 UTC 2025
# Append one new data row (timestamp + metric) after the existing data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A65").Value = "2025-04-29 09:10:15"
$ws.Range("B65").Value = 190
